$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.010799999999993
$ws.Range("D4").Value = -7.005399999999996
$ws.Range("D6").Value = -8.378199999999994
$ws.Range("A9").Value = -20.52049999999998
$ws.Range("D10").Value = -7.628199999999993
$ws.Range("B11").Value = 5.629399999999999
$ws.Range("D11").Value = -8.235399999999997
$ws.Range("A18").Value = -22.95420000000003
$ws.Range("A20").Value = -22.11400000000002
$ws.Range("C21").Value = -13.1994
$ws.Range("D21").Value = -8.263999999999999
